$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "'3"
$ws.Range("B2").Value = "3 : résultats postés ou publiés après les 36 mois"
$ws.Range("C2").Value = "NCT00551551"
$ws.Range("F2").Value = "'2011"
$ws.Range("G2").Value = "Urinary Postpartum Handicap Prevention: Pelvic Floor Exercises vs Control. Multicentric Randomized Trial"
$ws.Range("I2").Value = "OTHER"

# Row 3
$ws.Range("A3").Value = "'3"
$ws.Range("B3").Value = "3 : résultats postés ou publiés après les 36 mois"
$ws.Range("C3").Value = "NCT01464073"
$ws.Range("F3").Value = "'2012"
$ws.Range("G3").Value = "Evaluation of the Effectiveness of Individualized Exercise Program, Combined With a Balanced Diet Rich in Fruits and Vegetables, on the Evolution of Body Fat in Overweight or Obese Women, as Part of Type 2 Diabetes Prevention in Reunion Island."
$ws.Range("H3").Value = "LIPOXmax-RUN"
$ws.Range("I3").Value = "DIETARY_SUPPLEMENT"

# Row 4
$ws.Range("A4").Value = "'4"
$ws.Range("B4").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("C4").Value = "NCT01794377"
$ws.Range("F4").Value = "'2013"
$ws.Range("G4").Value = "Influence of Muscle Strength on the Energy Cost of Walking in Obese Subjects"
$ws.Range("H4").Value = "OBELIX"
$ws.Range("I4").Value = "BEHAVIORAL"

# Row 5
$ws.Range("A5").Value = "'4"
$ws.Range("B5").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("C5").Value = "NCT02275429"
$ws.Range("F5").Value = "'2014"
$ws.Range("G5").Value = "Etude Des Troubles métaboliques Induits Par Une Course à Pied très Longue Distance, `" La Diagonale Des Fous `" de La Réunion / Metabolic Disorders in Ultramarathon Runners of the Madmen's Diagonal Race on Reunion Island"
$ws.Range("H5").Value = "METARUN"
$ws.Range("I5").Value = "PROCEDURE"

# Row 6
$ws.Range("A6").Value = "'3"
$ws.Range("B6").Value = "3 : résultats postés ou publiés après les 36 mois"
$ws.Range("C6").Value = "NCT02000674"
$ws.Range("D6").Value = "2013-001438-16"
$ws.Range("F6").Value = "'2016"
$ws.Range("G6").Value = "Succinylcholine vs Rocuronium for Prehospital Emergency Intubation : a Randomized Trial"
$ws.Range("H6").Value = "CURASMUR"
$ws.Range("I6").Value = "DRUG"

# Row 7
$ws.Range("A7").Value = "'3"
$ws.Range("B7").Value = "3 : résultats postés ou publiés après les 36 mois"
$ws.Range("C7").Value = "NCT01425866"
$ws.Range("F7").Value = "'2016"
$ws.Range("G7").Value = "Multicenter Randomized Trial of Structured Educational Intervention at the Community Level in Insufficiently Controlled Patients With Type 2 Diabetes in Reunion Island"
$ws.Range("H7").Value = "ERMIES"
$ws.Range("I7").Value = "BEHAVIORAL"

# Row 8
$ws.Range("A8").Value = "'4"
$ws.Range("B8").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("C8").Value = "NCT02900989"
$ws.Range("F8").Value = "'2017"
$ws.Range("G8").Value = "Etude de Validation de la Version française d'un Questionnaire `" Ask Suicide-Screening Questions `" (ASQ) Dans Une Population de Patients Adolescents Pris en Charge en unité d'Urgences pédiatriques"
$ws.Range("H8").Value = "ASQ-Fr"
$ws.Range("I8").Value = "BEHAVIORAL"

# Row 9
$ws.Range("A9").Value = "'1"
$ws.Range("B9").Value = "1 : résultats postés ou publiés dans les 12 mois"
$ws.Range("C9").Value = "NCT03226834"
$ws.Range("F9").Value = "'2018"
$ws.Range("G9").Value = "Comparison of Musicotherapy Sessions Using Patient Play-list Versus U Sequence Music Care Medical Device on Pre-operative Anxiety in Women Undergoing Gynecological Surgery: a Randomized-control Study"
$ws.Range("H9").Value = "MUANX"
$ws.Range("I9").Value = "DEVICE"

# Row 10
$ws.Range("A10").Value = "'4"
$ws.Range("B10").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("C10").Value = "NCT03271112"
$ws.Range("F10").Value = "'2019"
$ws.Range("G10").Value = "Frailty Prevention in Elderly People From Reunion Island: Effects of Adaptated Exercises on Physical Performance"
$ws.Range("H10").Value = "5P-PILOTE"
$ws.Range("I10").Value = "OTHER"

# Row 11
$ws.Range("A11").Value = "'2"
$ws.Range("B11").Value = "2 : résultats postés ou publiés entre 12 et 36 mois"
$ws.Range("C11").Value = "NCT01537601"
$ws.Range("F11").Value = "'2019"
$ws.Range("G11").Value = "Effect of Circumcision on the Risk of Febrile Urinary Tract Infections in Children With Posterior Urethral Valves."
$ws.Range("H11").Value = "CIRCUP"
$ws.Range("I11").Value = "PROCEDURE"

# Row 12
$ws.Range("A12").Value = "'4"
$ws.Range("B12").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("C12").Value = "NCT04909411"
$ws.Range("F12").Value = "'2021"
$ws.Range("G12").Value = "Consequences of a Maternal-fetal Chikungunya Virus Infection. Neurocognitive and Sensory Assessment Around the Age of 13."
$ws.Range("H12").Value = "CHIK13+"
$ws.Range("I12").Value = "OTHER"

# Row 13
$ws.Range("A13").Value = "'4"
$ws.Range("B13").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("C13").Value = "NCT04768621"
$ws.Range("F13").Value = "'2021"
$ws.Range("G13").Value = "Health Consequences of Wintering in the French Southern and Antarctic Territories"
$ws.Range("H13").Value = "MediTAAF"
$ws.Range("I13").Value = "OTHER"

# Row 14
$ws.Range("A14").Value = "'1"
$ws.Range("B14").Value = "1 : résultats postés ou publiés dans les 12 mois"
$ws.Range("C14").Value = "NCT04459221"
$ws.Range("D14").Value = "2020-002332-73"
$ws.Range("F14").Value = "'2021"
$ws.Range("G14").Value = "Study of the Impact of a School Program Combining - Promotion of HPV Vaccination and HPV Vaccine Offer in Middle School - on Adherence to HPV Vaccination in Middle School Students"
$ws.Range("H14").Value = "PROM SSCOL"
$ws.Range("I14").Value = "DRUG"

# Row 15
$ws.Range("A15").Value = "'4"
$ws.Range("B15").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("C15").Value = "NCT05098925"
$ws.Range("F15").Value = "'2021"
$ws.Range("G15").Value = "Study of Thermoregulatory Processes in Ultra-endurance Runners in a Hot and Humid Environment"
$ws.Range("H15").Value = "ERUPTION-2"
$ws.Range("I15").Value = "OTHER"

# Row 16
$ws.Range("A16").Value = "'4"
$ws.Range("B16").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("C16").Value = "NCT05413720"
$ws.Range("F16").Value = "'2022"
$ws.Range("G16").Value = "Study of the Innate Immune Response to the Acute Phase of Human Leptospirosis - IMMUNOLEPTO"
$ws.Range("H16").Value = "IMMUNOLEPTO"
$ws.Range("I16").Value = "OTHER"

# Row 17
$ws.Range("A17").Value = "'4"
$ws.Range("B17").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("C17").Value = "NCT05237180"
$ws.Range("F17").Value = "'2022"
$ws.Range("G17").Value = "Evaluation of the Effect of a Spatial Localization Training Program on Auditory Comprehension in Noise in Bi-implanted Subjects With Post-lingual Deafness"
$ws.Range("H17").Value = "CAudiBruit"
$ws.Range("I17").Value = "OTHER"

# Row 18
$ws.Range("A18").Value = "'4"
$ws.Range("B18").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("C18").Value = "NCT05424913"
$ws.Range("F18").Value = "'2023"
$ws.Range("G18").Value = "Exploratory Study of the Relationships Between the Biomarkers of Inflammation, Lipidome and Insulin Resistance and Disorders of Glycemic Regulation in a Cohort of Insulin-resistant Subjects Due to Excess Weight or Dunnigan's Lipodystrophy"
$ws.Range("H18").Value = "IRAP-DUN 2"
$ws.Range("I18").Value = "OTHER"

# Row 19
$ws.Range("A19").Value = "'4"
$ws.Range("B19").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("C19").Value = "NCT05367180"
$ws.Range("F19").Value = "'2023"
$ws.Range("G19").Value = "Impact of a Prevention Program on Sun Risks in Primary School in Tropical French Region"
$ws.Range("H19").Value = "PRESOLRE"
$ws.Range("I19").Value = "OTHER"

# Row 20
$ws.Range("A20").Value = "'4"
$ws.Range("B20").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("C20").Value = "NCT05598138"
$ws.Range("F20").Value = "'2023"
$ws.Range("G20").Value = "Clinical and Biological Strokes Collection in Reunion Island"
$ws.Range("H20").Value = "CoBRA"
$ws.Range("I20").Value = "BIOLOGICAL"

# Row 21
$ws.Range("A21").Value = "'4"
$ws.Range("B21").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("C21").Value = "NCT06471335"
$ws.Range("F21").Value = "'2024"
$ws.Range("G21").Value = "Fetal Alcohol Spectrum Disorder (FASD): Clinical Description and Search for Epigenetic Biomarkers for Diagnostic Purposes."
$ws.Range("H21").Value = "EPI-TSAF"
$ws.Range("I21").Value = "DIAGNOSTIC_TEST"

# Force text-number cells back to default (Normal) style so no stray quotePrefix style lingers visually
$ws.Range("A2").Style = "Normal"
$ws.Range("F2").Style = "Normal"
$ws.Range("A3").Style = "Normal"
$ws.Range("F3").Style = "Normal"
$ws.Range("A4").Style = "Normal"
$ws.Range("F4").Style = "Normal"
$ws.Range("A5").Style = "Normal"
$ws.Range("F5").Style = "Normal"
$ws.Range("A6").Style = "Normal"
$ws.Range("F6").Style = "Normal"
$ws.Range("A7").Style = "Normal"
$ws.Range("F7").Style = "Normal"
$ws.Range("A8").Style = "Normal"
$ws.Range("F8").Style = "Normal"
$ws.Range("A9").Style = "Normal"
$ws.Range("F9").Style = "Normal"
$ws.Range("A10").Style = "Normal"
$ws.Range("F10").Style = "Normal"
$ws.Range("A11").Style = "Normal"
$ws.Range("F11").Style = "Normal"
$ws.Range("A12").Style = "Normal"
$ws.Range("F12").Style = "Normal"
$ws.Range("A13").Style = "Normal"
$ws.Range("F13").Style = "Normal"
$ws.Range("A14").Style = "Normal"
$ws.Range("F14").Style = "Normal"
$ws.Range("A15").Style = "Normal"
$ws.Range("F15").Style = "Normal"
$ws.Range("A16").Style = "Normal"
$ws.Range("F16").Style = "Normal"
$ws.Range("A17").Style = "Normal"
$ws.Range("F17").Style = "Normal"
$ws.Range("A18").Style = "Normal"
$ws.Range("F18").Style = "Normal"
$ws.Range("A19").Style = "Normal"
$ws.Range("F19").Style = "Normal"
$ws.Range("A20").Style = "Normal"
$ws.Range("F20").Style = "Normal"
$ws.Range("A21").Style = "Normal"
$ws.Range("F21").Style = "Normal"